$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: 'Datos actualizados a 1 de Agosto de 2020 a las 23:14' -> 'Datos actualizados a 2 de Agosto de 2020 a las 00:31'
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 00:31"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 4758204
$ws.Range("C4").Value = 52315
$ws.Range("D4").Value = 2354419
$ws.Range("E4").Value = 2246048
$ws.Range("G4").Value = 990
$ws.Range("H4").Value = 157737

# Row 14: 'Reino Unido' -> 'Colombia'
$ws.Range("A14").Value = "Colombia"
$ws.Range("B14").Value = 306181
$ws.Range("C14").Value = 10673
$ws.Range("D14").Value = 160708
$ws.Range("E14").Value = 135143
$ws.Range("G14").Value = 225
$ws.Range("H14").Value = 10330

# Row 15: 'Colombia' -> 'Reino Unido'
$ws.Range("A15").Value = "Reino Unido"
$ws.Range("B15").Value = 303952
$ws.Range("C15").Value = 771
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 74
$ws.Range("H15").Value = 46193

# Row 21: 'Alemania' -> 'Alemania'
$ws.Range("B21").Value = 211077
$ws.Range("C21").Value = 412
$ws.Range("E21").Value = 8251

# Row 25: 'Canada' -> 'Canada'
$ws.Range("B25").Value = 116599
$ws.Range("C25").Value = 287
$ws.Range("D25").Value = 101436
$ws.Range("E25").Value = 6222

# Row 52: 'Barein' -> 'Barein'
$ws.Range("B52").Value = 41190
$ws.Range("C52").Value = 208
$ws.Range("D52").Value = 38211
$ws.Range("E52").Value = 2832

# Row 81: 'Sudan' -> 'Bulgaria'
$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 11836
$ws.Range("C81").Value = 146
$ws.Range("D81").Value = 6396
$ws.Range("E81").Value = 5055
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 385

# Row 82: 'Bulgaria' -> 'Sudan'
$ws.Range("A82").Value = "Sudan"
$ws.Range("B82").Value = 11738
$ws.Range("C82").Value = 94
$ws.Range("D82").Value = 6137
$ws.Range("E82").Value = 4849
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 752

# Row 90: 'Tayikistan' -> 'Gabon'
$ws.Range("A90").Value = "Gabon"
$ws.Range("B90").Value = 7531
$ws.Range("C90").Value = 179
$ws.Range("D90").Value = 5223
$ws.Range("E90").Value = 2258
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 50

# Row 91: 'Finlandia' -> 'Tayikistan'
$ws.Range("A91").Value = "Tayikistan"
$ws.Range("B91").Value = 7451
$ws.Range("C91").Value = 42
$ws.Range("D91").Value = 6233
$ws.Range("E91").Value = 1158
$ws.Range("H91").Value = 60

# Row 92: 'Haiti' -> 'Finlandia'
$ws.Range("A92").Value = "Finlandia"
$ws.Range("B92").Value = 7443
$ws.Range("C92").Value = 11
$ws.Range("D92").Value = 6950
$ws.Range("E92").Value = 164
$ws.Range("H92").Value = 329

# Row 93: 'Gabon' -> 'Haiti'
$ws.Range("A93").Value = "Haiti"
$ws.Range("B93").Value = 7424
$ws.Range("C93").Value = 12
$ws.Range("D93").Value = 4606
$ws.Range("E93").Value = 2657
$ws.Range("H93").Value = 161

# Row 108: 'Maldivas' -> 'Maldivas'
$ws.Range("E108").Value = 1319
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 17

# Row 110: 'Libia' -> 'Zimbabue'
$ws.Range("A110").Value = "Zimbabue"
$ws.Range("B110").Value = 3659
$ws.Range("C110").Value = 490
$ws.Range("D110").Value = 1011
$ws.Range("E110").Value = 2579
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 69

# Row 111: 'Hong Kong' -> 'Libia'
$ws.Range("A111").Value = "Libia"
$ws.Range("B111").Value = 3621
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 618
$ws.Range("E111").Value = 2929
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 74

# Row 112: 'Tailandia' -> 'Hong Kong'
$ws.Range("A112").Value = "Hong Kong"
$ws.Range("B112").Value = 3398
$ws.Range("C112").Value = 125
$ws.Range("D112").Value = 1858
$ws.Range("E112").Value = 1507
$ws.Range("G112").Value = 6
$ws.Range("H112").Value = 33

# Row 113: 'Somalia' -> 'Tailandia'
$ws.Range("A113").Value = "Tailandia"
$ws.Range("B113").Value = 3312
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 3135
$ws.Range("E113").Value = 119
$ws.Range("H113").Value = 58

# Row 114: 'Congo' -> 'Somalia'
$ws.Range("A114").Value = "Somalia"
$ws.Range("B114").Value = 3212
$ws.Range("D114").Value = 1562
$ws.Range("E114").Value = 1557
$ws.Range("H114").Value = 93

# Row 115: 'Montenegro' -> 'Congo'
$ws.Range("A115").Value = "Congo"
$ws.Range("B115").Value = 3200
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 829
$ws.Range("E115").Value = 2317
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 54

# Row 116: 'Zimbabue' -> 'Montenegro'
$ws.Range("A116").Value = "Montenegro"
$ws.Range("B116").Value = 3198
$ws.Range("C116").Value = 86
$ws.Range("D116").Value = 1293
$ws.Range("E116").Value = 1855
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 50

# Row 129: 'Ruanda' -> 'Ruanda'
$ws.Range("B129").Value = 2042
$ws.Range("C129").Value = 20
$ws.Range("D129").Value = 1119
$ws.Range("E129").Value = 918
